$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must be forced to Text format
# first, otherwise Excel auto-converts the literal into a Number cell (losing
# formatting like trailing zeros, e.g. "1.00" -> 1).
$textCells = @("D5","D6","D8","D10","D15","D16","D20","D25","D28","D29","D30","D37","D40","D41","D43","D44","D46","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell
$ws.Range('D2').Value = '26.635.65'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '1.596.74'
$ws.Range('E3').Value = '  -0.22%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '210.67'
$ws.Range('D6').Value = '0.510'
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').Value = '0.0615'
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('E9').Value = '  -0.38%  '
$ws.Range('D10').Value = '19.59'
$ws.Range('E10').Value = '  +0.16%  '
$ws.Range('E11').Value = '  +0.14%  '
$ws.Range('D12').Value = '1.821.13'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').Value = '1.598.98'
$ws.Range('E13').Value = '  +1.45%  '
$ws.Range('E14').Value = '  -0.20%  '
$ws.Range('D15').Value = '0.523'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '64.61'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '26.607.87'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('E18').Value = '  -2.61%  '
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '208.63'
$ws.Range('E20').Value = '  -0.36%  '
$ws.Range('E21').Value = '  -1.55%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('E23').Value = '  -3.20%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '144.06'
$ws.Range('E25').Value = '  +0.76%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('D28').Value = '0.114'
$ws.Range('E28').Value = '  -0.94%  '
$ws.Range('D29').Value = '15.25'
$ws.Range('E29').Value = '  -0.68%  '
$ws.Range('D30').Value = '0.0506'
$ws.Range('E30').Value = '  -2.30%  '
$ws.Range('E31').Value = '  -0.43%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('D34').Value = '1.279.58'
$ws.Range('E34').Value = '  -1.10%  '
$ws.Range('E35').Value = '  +0.55%  '
$ws.Range('E36').Value = '  +13.38%  '
$ws.Range('D37').Value = '0.601'
$ws.Range('E37').Value = '  -3.15%  '
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('D40').Value = '0.824'
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('D41').Value = '5.43'
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('D43').Value = '0.771'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('D44').Value = '62.85'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('D45').Value = '1.732.82'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').Value = '89.46'
$ws.Range('E46').Value = '  -1.84%  '
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('E48').Value = '  +2.31%  '
$ws.Range('D49').Value = '0.0513'
$ws.Range('E49').Value = '  +0.72%  '
$ws.Range('B50').Value = 'USDD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  -0.05%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.45'
$ws.Range('E51').Value = '  +1.19%  '
